# Apply attendance updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where Total Attendance Count (D) and Real (E) become 1
$realRows = @(4, 5, 6, 9, 10, 11, 12, 13, 14, 15, 18)
foreach ($r in $realRows) {
    $ws.Cells.Item($r, 4).Value = 1   # column D
    $ws.Cells.Item($r, 5).Value = 1   # column E
}

# Rows where Absent (H) becomes 1
$absentRows = @(3, 7, 8, 16, 17)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1   # column H
}

# Row 3 also has Invalid (G) become 1
$ws.Cells.Item(3, 7).Value = 1   # column G
